$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels in row 1 to uppercase language suffixes
$ws.Range("A1").Value = "instr_msg_EN"
$ws.Range("C1").Value = "instr_msg_ES"
$ws.Range("D1").Value = "instr_msg_FR"

# Swap the images_sizeW / images_sizeH values on row 2
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.5

# Update the active selection to F2, matching the saved selection state
$ws.Range("F2").Select()
